# Weekly fruit/hortaliza data update.
# A new weekly price-report row is inserted at row 187 (pushing the existing
# rows 187-209 down to 188-210), and the new row is populated with this
# week's "Ajo" (Chino / Primera) price data for the Femacal de La Calera
# market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 187; everything from 187 downward shifts to 188+.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new observation.
$ws.Cells.Item(187, 1).Value = 3
$ws.Cells.Item(187, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(187, 3).Value = "Coquimbo"
$ws.Cells.Item(187, 4).Value = 44449
$ws.Cells.Item(187, 5).Value = 5
$ws.Cells.Item(187, 6).Value = 100112003
$ws.Cells.Item(187, 7).Value = "Ajo"
$ws.Cells.Item(187, 8).Value = "Chino"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 93
$ws.Cells.Item(187, 11).Value = 15000
$ws.Cells.Item(187, 12).Value = 15500
$ws.Cells.Item(187, 13).Value = 15258
$ws.Cells.Item(187, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(187, 15).Value = "China"
$ws.Cells.Item(187, 16).Value = 1526
$ws.Cells.Item(187, 17).Value = 10
$ws.Cells.Item(187, 18).Value = "Hortaliza"
